$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2084.5386
$ws.Range("I106").Value = 1511
$ws.Range("K106").Value = 1511
$ws.Range("M106").Value = -880

$ws.Range("H129").Value = 937.7027
$ws.Range("I129").Value = 384.07144
$ws.Range("J129").Value = 1274.6957
$ws.Range("K129").Value = 1152.21432
$ws.Range("L129").Value = 3824.0871
$ws.Range("M129").Value = 3847.78568
$ws.Range("N129").Value = -13824.0871

$ws.Range("H137").Value = 1257.7333
$ws.Range("I137").Value = 1256.6538
$ws.Range("J137").Value = 1259.2106
$ws.Range("K137").Value = 3769.9614
$ws.Range("L137").Value = 3777.6318
$ws.Range("M137").Value = -1219.9614
$ws.Range("N137").Value = -8877.631799999999

$ws.Range("H138").Value = 16396540
$ws.Range("I138").Value = 1385.725
$ws.Range("J138").Value = 47625404
$ws.Range("K138").Value = 4157.174999999999
$ws.Range("L138").Value = 142876212
$ws.Range("M138").Value = 982.8250000000007
$ws.Range("N138").Value = -142886492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4680.1167
$ws.Range("I32").Value = 2977.9075
$ws.Range("K32").Value = 2977.9075
$ws.Range("M32").Value = -2690.9075

$ws.Range("H132").Value = 3091.7
$ws.Range("I132").Value = 1411.4546
$ws.Range("J132").Value = 4064.4736
$ws.Range("K132").Value = 4234.3638
$ws.Range("L132").Value = 12193.4208
$ws.Range("M132").Value = -1704.3638
$ws.Range("N132").Value = -17253.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1914.2858
$ws.Range("I105").Value = 2333.3333
$ws.Range("K105").Value = 2333.3333
$ws.Range("M105").Value = -586.3332999999998

$ws.Range("H134").Value = 2494.8147
$ws.Range("I134").Value = 1703.2354
$ws.Range("J134").Value = 3840.5
$ws.Range("K134").Value = 5109.706200000001
$ws.Range("L134").Value = 11521.5
$ws.Range("M134").Value = -2574.706200000001
$ws.Range("N134").Value = -16591.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2226.0625
$ws.Range("I132").Value = 1748.88
$ws.Range("J132").Value = 3930.2856
$ws.Range("K132").Value = 5246.64
$ws.Range("L132").Value = 11790.8568
$ws.Range("M132").Value = -2716.64
$ws.Range("N132").Value = -16850.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 651350.4399999999
$ws.Range("I5").Value = 629.625
$ws.Range("J5").Value = 925338.2
$ws.Range("K5").Value = 1888.875
$ws.Range("L5").Value = 2776014.6
$ws.Range("M5").Value = -1776.875
$ws.Range("N5").Value = -2776238.6

$ws.Range("H68").Value = 992.4
$ws.Range("I68").Value = 774.1627999999999
$ws.Range("J68").Value = 1157.035
$ws.Range("K68").Value = 2322.4884
$ws.Range("L68").Value = 3471.105
$ws.Range("M68").Value = -1511.4884
$ws.Range("N68").Value = -5093.105

$ws.Range("H71").Value = 992.4
$ws.Range("I71").Value = 774.1627999999999
$ws.Range("J71").Value = 1157.035
$ws.Range("K71").Value = 6967.4652
$ws.Range("L71").Value = 10413.315
$ws.Range("M71").Value = -2911.4652
$ws.Range("N71").Value = -18525.315

$ws.Range("H107").Value = 720.54
$ws.Range("I107").Value = 658.75
$ws.Range("J107").Value = 799.1818
$ws.Range("K107").Value = 1976.25
$ws.Range("L107").Value = 2397.5454
$ws.Range("M107").Value = -56.25
$ws.Range("N107").Value = -6237.5454

$ws.Range("H109").Value = 3220.2
$ws.Range("I109").Value = 867.6667
$ws.Range("J109").Value = 6749
$ws.Range("K109").Value = 2603.0001
$ws.Range("L109").Value = 20247
$ws.Range("M109").Value = -1563.0001
$ws.Range("N109").Value = -22327

$ws.Range("H110").Value = 4833.3335
$ws.Range("I110").Value = 500
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 1500
$ws.Range("L110").Value = 21000
$ws.Range("M110").Value = 2590
$ws.Range("N110").Value = -29180

$ws.Range("H112").Value = 1630
$ws.Range("J112").Value = 1687.5
$ws.Range("L112").Value = 5062.5
$ws.Range("N112").Value = -7278.5

$ws.Range("H122").Value = 494
$ws.Range("I122").Value = 353.36
$ws.Range("K122").Value = 3180.24
$ws.Range("M122").Value = -730.2400000000002

$ws.Range("H131").Value = 5106.1113
$ws.Range("I131").Value = 624.875
$ws.Range("J131").Value = 6992.9473
$ws.Range("K131").Value = 1874.625
$ws.Range("L131").Value = 20978.8419
$ws.Range("M131").Value = 3165.375
$ws.Range("N131").Value = -31058.8419

$ws.Range("H132").Value = 1418.6666
$ws.Range("I132").Value = 1604.5
$ws.Range("J132").Value = 1270
$ws.Range("K132").Value = 14440.5
$ws.Range("L132").Value = 11430
$ws.Range("M132").Value = -11910.5
$ws.Range("N132").Value = -16490

$ws.Range("H135").Value = 651350.4399999999
$ws.Range("I135").Value = 629.625
$ws.Range("J135").Value = 925338.2
$ws.Range("K135").Value = 5666.625
$ws.Range("L135").Value = 8328043.8
$ws.Range("M135").Value = -3131.625
$ws.Range("N135").Value = -8333113.8

$ws.Range("H137").Value = 2480.8928
$ws.Range("I137").Value = 1858.2778
$ws.Range("J137").Value = 3601.6
$ws.Range("K137").Value = 5574.8334
$ws.Range("L137").Value = 10804.8
$ws.Range("M137").Value = -474.8334000000004
$ws.Range("N137").Value = -21004.8

$ws.Range("H140").Value = 1732.8214
$ws.Range("I140").Value = 1189.56
$ws.Range("J140").Value = 6260
$ws.Range("K140").Value = 3568.68
$ws.Range("L140").Value = 18780
$ws.Range("M140").Value = 1611.32
$ws.Range("N140").Value = -29140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 55000
$ws.Range("J47").Value = 55000
$ws.Range("L47").Value = 55000
$ws.Range("N47").Value = -56136

$ws.Range("H48").Value = 33465
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H55").Value = 4777.6665
$ws.Range("J55").Value = 6166.5
$ws.Range("L55").Value = 6166.5
$ws.Range("N55").Value = -6820.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 14500
$ws.Range("J3").Value = 14500
$ws.Range("L3").Value = 14500
$ws.Range("N3").Value = -14724

$ws.Range("H14").Value = 16500
$ws.Range("J14").Value = 16500
$ws.Range("L14").Value = 16500
$ws.Range("N14").Value = -16844

$ws.Range("H15").Value = 14500
$ws.Range("J15").Value = 14500
$ws.Range("L15").Value = 14500
$ws.Range("N15").Value = -14840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 911.1111
$ws.Range("J14").Value = 1800
$ws.Range("L14").Value = 1800
$ws.Range("N14").Value = -2136

$ws.Range("H100").Value = 67456.53
$ws.Range("I100").Value = 772.44446
$ws.Range("J100").Value = 167482.67
$ws.Range("K100").Value = 1544.88892
$ws.Range("L100").Value = 334965.34
$ws.Range("M100").Value = -1003.88892
$ws.Range("N100").Value = -336047.34

$ws.Range("H132").Value = 2799.2974
$ws.Range("I132").Value = 2152.9375
$ws.Range("J132").Value = 3291.762
$ws.Range("K132").Value = 6458.8125
$ws.Range("L132").Value = 9875.286
$ws.Range("M132").Value = -3928.8125
$ws.Range("N132").Value = -14935.286
